$d = $word.ActiveDocument

# Replace "Apples" with "Apples<newline>Pears<newline>Grapes" across the
# whole document. Soft line breaks (Shift+Enter) are represented in
# PowerShell/VBA text as Chr(11) (vertical tab, written here as the
# backtick-v escape) and are written to the OOXML as <w:br/> elements
# rather than new paragraphs, which is how replace_docx_text now
# preserves embedded newlines in the replacement text.
$d.Content.Find.Execute("Apples", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Apples`vPears`vGrapes", 2)
